$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 165.89473
$ws.Range("I2").Value = 130.66667
$ws.Range("K2").Value = 130.66667
$ws.Range("M2").Value = -17.66667000000001
$ws.Range("H12").Value = 259.64285
$ws.Range("I12").Value = 255.5
$ws.Range("J12").Value = 284.5
$ws.Range("K12").Value = 255.5
$ws.Range("L12").Value = 284.5
$ws.Range("M12").Value = -85.5
$ws.Range("N12").Value = -624.5
$ws.Range("H18").Value = 3453
$ws.Range("I18").Value = 1296.7142
$ws.Range("J18").Value = 11000
$ws.Range("K18").Value = 1296.7142
$ws.Range("L18").Value = 11000
$ws.Range("M18").Value = -1012.7142
$ws.Range("N18").Value = -11568
$ws.Range("H19").Value = 29799.2
$ws.Range("I19").Value = 4999
$ws.Range("J19").Value = 40427.855
$ws.Range("K19").Value = 4999
$ws.Range("L19").Value = 40427.855
$ws.Range("M19").Value = -4824
$ws.Range("N19").Value = -40777.855
$ws.Range("H28").Value = 29989.941
$ws.Range("I28").Value = 33915.766
$ws.Range("J28").Value = 546.25
$ws.Range("K28").Value = 33915.766
$ws.Range("L28").Value = 546.25
$ws.Range("M28").Value = -33430.766
$ws.Range("N28").Value = -1516.25
$ws.Range("H34").Value = 3663.6667
$ws.Range("I34").Value = 3663.6667
$ws.Range("K34").Value = 3663.6667
$ws.Range("M34").Value = -3460.6667
$ws.Range("H36").Value = 3663.6667
$ws.Range("I36").Value = 3663.6667
$ws.Range("K36").Value = 3663.6667
$ws.Range("M36").Value = -2948.6667
$ws.Range("H40").Value = 2090.4546
$ws.Range("I40").Value = 997.5
$ws.Range("K40").Value = 997.5
$ws.Range("M40").Value = -822.5
$ws.Range("H55").Value = 99.666664
$ws.Range("I55").Value = 99.666664
$ws.Range("K55").Value = 99.666664
$ws.Range("M55").Value = 114.333336
$ws.Range("H98").Value = 2426.1462
$ws.Range("I98").Value = 2282.1765
$ws.Range("K98").Value = 2282.1765
$ws.Range("M98").Value = -784.1765
$ws.Range("H111").Value = 31294.545
$ws.Range("I111").Value = 1892.4286
$ws.Range("J111").Value = 82748.25
$ws.Range("K111").Value = 5677.2858
$ws.Range("L111").Value = 248244.75
$ws.Range("M111").Value = -2610.2858
$ws.Range("N111").Value = -254378.75
$ws.Range("H116").Value = 20157882
$ws.Range("I116").Value = 14775567
$ws.Range("J116").Value = 27782828
$ws.Range("K116").Value = 14775567
$ws.Range("L116").Value = 27782828
$ws.Range("M116").Value = -14772125
$ws.Range("N116").Value = -27789712
$ws.Range("H122").Value = 2426.1462
$ws.Range("I122").Value = 2282.1765
$ws.Range("K122").Value = 6846.529500000001
$ws.Range("M122").Value = -4396.529500000001
$ws.Range("H132").Value = 2721.75
$ws.Range("I132").Value = 2392.7307
$ws.Range("K132").Value = 7178.1921
$ws.Range("M132").Value = -4648.1921
$ws.Range("H135").Value = 901.6667
$ws.Range("I135").Value = 901.6667
$ws.Range("K135").Value = 8115.0003
$ws.Range("M135").Value = -5580.0003
$ws.Range("H137").Value = 3421.1724
$ws.Range("I137").Value = 2809.7222
$ws.Range("J137").Value = 4421.727
$ws.Range("K137").Value = 8429.1666
$ws.Range("L137").Value = 13265.181
$ws.Range("M137").Value = -5879.1666
$ws.Range("N137").Value = -18365.181

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 730.0625
$ws.Range("I2").Value = 932.2
$ws.Range("K2").Value = 932.2
$ws.Range("M2").Value = -819.2
$ws.Range("H32").Value = 4701.877
$ws.Range("I32").Value = 3973.754
$ws.Range("K32").Value = 3973.754
$ws.Range("M32").Value = -3686.754
$ws.Range("H45").Value = 3201.625
$ws.Range("I45").Value = 3099.8
$ws.Range("K45").Value = 3099.8
$ws.Range("M45").Value = -2722.8
$ws.Range("H61").Value = 3381.818
$ws.Range("I61").Value = 1570.6154
$ws.Range("K61").Value = 1570.6154
$ws.Range("M61").Value = -1358.6154
$ws.Range("H74").Value = 2479.2
$ws.Range("I74").Value = 1711.6875
$ws.Range("K74").Value = 1711.6875
$ws.Range("M74").Value = -837.6875
$ws.Range("H77").Value = 2479.2
$ws.Range("I77").Value = 1711.6875
$ws.Range("K77").Value = 8558.4375
$ws.Range("M77").Value = -4190.4375
$ws.Range("H116").Value = 730.0625
$ws.Range("I116").Value = 932.2
$ws.Range("K116").Value = 932.2
$ws.Range("M116").Value = 1361.8
$ws.Range("H136").Value = 3381.818
$ws.Range("I136").Value = 1570.6154
$ws.Range("K136").Value = 4711.8462
$ws.Range("M136").Value = -2161.8462

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 730.0625
$ws.Range("I3").Value = 932.2
$ws.Range("K3").Value = 932.2
$ws.Range("M3").Value = -818.2
$ws.Range("H108").Value = 79999
$ws.Range("I108").Value = 79999
$ws.Range("K108").Value = 79999
$ws.Range("M108").Value = -76159
$ws.Range("H134").Value = 3979.4167
$ws.Range("J134").Value = 5444
$ws.Range("L134").Value = 16332
$ws.Range("N134").Value = -21402

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4456.5
$ws.Range("I99").Value = 4762.25
$ws.Range("J99").Value = 3997.875
$ws.Range("K99").Value = 4762.25
$ws.Range("L99").Value = 3997.875
$ws.Range("M99").Value = -3264.25
$ws.Range("N99").Value = -6993.875
$ws.Range("H105").Value = 912.2143
$ws.Range("I105").Value = 877.4
$ws.Range("J105").Value = 999.25
$ws.Range("K105").Value = 877.4
$ws.Range("L105").Value = 999.25
$ws.Range("M105").Value = 869.6
$ws.Range("N105").Value = -4493.25
$ws.Range("H107").Value = 601.6087
$ws.Range("I107").Value = 548.6842
$ws.Range("J107").Value = 853
$ws.Range("K107").Value = 548.6842
$ws.Range("L107").Value = 853
$ws.Range("M107").Value = 1371.3158
$ws.Range("N107").Value = -4693
$ws.Range("H126").Value = 4456.5
$ws.Range("I126").Value = 4762.25
$ws.Range("J126").Value = 3997.875
$ws.Range("K126").Value = 14286.75
$ws.Range("L126").Value = 11993.625
$ws.Range("M126").Value = -11816.75
$ws.Range("N126").Value = -16933.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 535
$ws.Range("I122").Value = 535
$ws.Range("K122").Value = 4815
$ws.Range("M122").Value = -2365
$ws.Range("H134").Value = 3378.2856
$ws.Range("I134").Value = 2274.6667
$ws.Range("K134").Value = 6824.000100000001
$ws.Range("M134").Value = -1754.000100000001
$ws.Range("H140").Value = 4316
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6226.909
$ws.Range("I80").Value = 3749.5
$ws.Range("K80").Value = 3749.5
$ws.Range("M80").Value = -2751.5
$ws.Range("H83").Value = 6226.909
$ws.Range("I83").Value = 3749.5
$ws.Range("K83").Value = 18747.5
$ws.Range("M83").Value = -13755.5
$ws.Range("H126").Value = 7547.478
$ws.Range("I126").Value = 10043
$ws.Range("J126").Value = 3665.5557
$ws.Range("K126").Value = 30129
$ws.Range("L126").Value = 10996.6671
$ws.Range("M126").Value = -27659
$ws.Range("N126").Value = -15936.6671
$ws.Range("H132").Value = 3036.9333
$ws.Range("I132").Value = 2798.742
$ws.Range("J132").Value = 3564.3572
$ws.Range("K132").Value = 8396.226000000001
$ws.Range("L132").Value = 10693.0716
$ws.Range("M132").Value = -5866.226000000001
$ws.Range("N132").Value = -15753.0716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21742454
$ws.Range("I7").Value = 35717012
$ws.Range("K7").Value = 35717012
$ws.Range("M7").Value = -35716900
$ws.Range("H40").Value = 940459.25
$ws.Range("I40").Value = 1502527.4
$ws.Range("K40").Value = 1502527.4
$ws.Range("M40").Value = -1502391.4
$ws.Range("H122").Value = 4204.4688
$ws.Range("I122").Value = 4003.158
$ws.Range("J122").Value = 4498.6924
$ws.Range("K122").Value = 12009.474
$ws.Range("L122").Value = 13496.0772
$ws.Range("M122").Value = -9559.474
$ws.Range("N122").Value = -18396.0772
$ws.Range("H126").Value = 21742454
$ws.Range("I126").Value = 35717012
$ws.Range("K126").Value = 107151036
$ws.Range("M126").Value = -107148566
$ws.Range("H136").Value = 4167.9624
$ws.Range("I136").Value = 3936.0222
$ws.Range("J136").Value = 5472.625
$ws.Range("K136").Value = 11808.0666
$ws.Range("L136").Value = 16417.875
$ws.Range("M136").Value = -9258.0666
$ws.Range("N136").Value = -21517.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 13999
$ws.Range("I26").Value = 13999
$ws.Range("K26").Value = 13999
$ws.Range("M26").Value = -13706
$ws.Range("H53").Value = 33999
$ws.Range("I53").Value = 33999
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 33999
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -33392
$ws.Range("N53").ClearContents()
$ws.Range("H107").Value = 214.5
$ws.Range("I107").Value = 148.33333
$ws.Range("K107").Value = 444.99999
$ws.Range("M107").Value = 1475.00001
$ws.Range("H126").Value = 1446.8334
$ws.Range("I126").Value = 1054.6471
$ws.Range("K126").Value = 3163.9413
$ws.Range("M126").Value = -693.9412999999995
$ws.Range("H132").Value = 342546.34
$ws.Range("I132").Value = 420064.28
$ws.Range("J132").Value = 4286.1816
$ws.Range("K132").Value = 1260192.84
$ws.Range("L132").Value = 12858.5448
$ws.Range("M132").Value = -1257662.84
$ws.Range("N132").Value = -17918.5448
